# "Generate Report for Handoff"
# The localization status report is being regenerated: the items that were
# previously handed back are now marked ready for a new handoff, and the
# handoff/generation timestamps are refreshed. Column widths on the
# "Status"-ish columns shrink because the new text ("Ready for handoff") is
# shorter than the old text ("Handed back: in sync with en-US").

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refreshed timestamps
$overview.Range("G2").Value = "2016-08-21 17:05:24"
$dede.Range("H2").Value = "2016-08-21 17:05:24"
$zhcn.Range("H2").Value = "2016-08-21 17:05:20"

# --- Column widths shrink to fit the new, shorter status text.
# "zh-cn"/"de-de" (Overview) and "Status" (zh-cn/de-de sheets) columns were
# sized for the long "Handed back: in sync with en-US" status string; now
# that the status text is the shorter "Ready for handoff", re-narrow them.
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
